$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.963.22'
$ws.Range("E2").Value = '  -3.27%  '
$ws.Range("D3").Value = '2.284.45'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.00'
$ws.Range("E5").Value = '  -2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.22'
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("D9").Value = '2.284.32'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("E10").Value = '  -3.29%  '
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.62'
$ws.Range("D14").Value = '2.675.20'
$ws.Range("E14").Value = '  -3.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.49'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").Value = '54.034.73'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("D18").Value = '2.268.34'
$ws.Range("E18").Value = '  -4.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.91'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '301.67'
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.41'
$ws.Range("E22").Value = '  +3.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.74'
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("D28").Value = '2.389.33'
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.26'
$ws.Range("E31").Value = '  -3.44%  '
$ws.Range("E32").Value = '  -2.89%  '
$ws.Range("D33").Value = '0.0₃0680'
$ws.Range("E33").Value = '  -3.18%  '
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.55'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.869'
$ws.Range("E40").Value = '  +4.42%  '
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.43'
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("E43").Value = '  +0.94%  '
$ws.Range("E44").Value = '  +0.73%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.38'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0886'
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.543'
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '237.11'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("E51").Value = '  +0.63%  '
